$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to retain an exact text value (avoids Excel
# auto-converting numeric-looking strings like "578.53" or "1.00" into
# numbers, which would drop formatting such as trailing zeros).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "69.607.83"
$ws.Range("E2").Value = "  -1.96%  "
Set-TextValue $ws.Range("D3") "3.564.68"
$ws.Range("E3").Value = "  -2.08%  "
$ws.Range("E4").Value = "  +0.15%  "
Set-TextValue $ws.Range("D5") "578.53"
$ws.Range("E5").Value = "  -2.70%  "
Set-TextValue $ws.Range("D6") "188.44"
$ws.Range("E6").Value = "  -2.01%  "
Set-TextValue $ws.Range("D7") "0.630"
$ws.Range("E7").Value = "  -3.85%  "
Set-TextValue $ws.Range("D8") "3.565.36"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("E9").Value = "  +0.08%  "
Set-TextValue $ws.Range("D10") "0.176"
$ws.Range("E10").Value = "  -2.28%  "
Set-TextValue $ws.Range("D11") "0.658"
$ws.Range("E11").Value = "  -1.20%  "
Set-TextValue $ws.Range("D12") "55.63"
$ws.Range("E12").Value = "  -4.31%  "
Set-TextValue $ws.Range("D13") "0.0000300"
$ws.Range("E13").Value = "  +1.14%  "
Set-TextValue $ws.Range("D14") "9.56"
$ws.Range("E14").Value = "  -2.29%  "
Set-TextValue $ws.Range("D15") "4.143.18"
$ws.Range("E15").Value = "  -1.74%  "
Set-TextValue $ws.Range("D16") "19.67"
$ws.Range("E16").Value = "  +1.14%  "
Set-TextValue $ws.Range("D17") "3.568.85"
$ws.Range("E17").Value = "  -1.81%  "
Set-TextValue $ws.Range("D18") "69.644.34"
$ws.Range("E18").Value = "  -1.61%  "
Set-TextValue $ws.Range("D19") "12.56"
$ws.Range("E19").Value = "  -0.73%  "
$ws.Range("E20").Value = "  -0.17%  "
Set-TextValue $ws.Range("D21") "1.03"
$ws.Range("E21").Value = "  -2.22%  "
Set-TextValue $ws.Range("D22") "473.46"
$ws.Range("E22").Value = "  -4.29%  "
Set-TextValue $ws.Range("D23") "19.22"
$ws.Range("E23").Value = "  +12.50%  "
Set-TextValue $ws.Range("D24") "5.01"
$ws.Range("E24").Value = "  -8.19%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D25") "4.36"
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "95.26"
$ws.Range("E26").Value = "  +4.39%  "
Set-TextValue $ws.Range("D27") "3.00"
$ws.Range("E27").Value = "  -4.21%  "
Set-TextValue $ws.Range("D28") "10.91"
$ws.Range("E28").Value = "  -3.57%  "
Set-TextValue $ws.Range("D29") "9.26"
$ws.Range("E29").Value = "  -2.55%  "
Set-TextValue $ws.Range("D30") "32.26"
$ws.Range("E30").Value = "  -0.59%  "
Set-TextValue $ws.Range("D31") "7.69"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("E32").Value = "  +0.22%  "
Set-TextValue $ws.Range("D33") "12.13"
$ws.Range("E33").Value = "  -1.16%  "
Set-TextValue $ws.Range("D34") "66.01"
$ws.Range("E34").Value = "  +0.91%  "
Set-TextValue $ws.Range("D35") "580.02"
$ws.Range("E35").Value = "  -5.94%  "
Set-TextValue $ws.Range("D36") "38.75"
$ws.Range("E36").Value = "  +1.20%  "
Set-TextValue $ws.Range("D37") "0.999"
$ws.Range("E37").Value = "  -0.07%  "
Set-TextValue $ws.Range("D38") "0.0₃0791"
$ws.Range("E38").Value = "  -5.43%  "
Set-TextValue $ws.Range("D39") "0.392"
$ws.Range("E39").Value = "  -3.94%  "
Set-TextValue $ws.Range("D40") "3.21"
$ws.Range("E40").Value = "  +16.41%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D41") "0.137"
$ws.Range("E41").Value = "  -7.66%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue $ws.Range("D42") "3.43"
$ws.Range("E42").Value = "  -7.44%  "
Set-TextValue $ws.Range("D43") "3.216.31"
$ws.Range("E43").Value = "  -3.84%  "
Set-TextValue $ws.Range("D44") "2.83"
$ws.Range("E44").Value = "  +3.93%  "
Set-TextValue $ws.Range("D45") "3.06"
$ws.Range("E45").Value = "  -0.69%  "
Set-TextValue $ws.Range("D46") "0.0440"
$ws.Range("E46").Value = "  -2.25%  "
Set-TextValue $ws.Range("D47") "3.36"
$ws.Range("E47").Value = "  -0.42%  "
Set-TextValue $ws.Range("D48") "9.40"
$ws.Range("E48").Value = "  +1.35%  "
$ws.Range("E49").Value = "  -1.08%  "
Set-TextValue $ws.Range("D50") "1.00"
$ws.Range("E50").Value = "  +0.07%  "
Set-TextValue $ws.Range("D51") "3.12"
$ws.Range("E51").Value = "  -6.83%  "
